{"js": "// Update the worksheet date and every \"a\u00f7b=\" problem cell to the new values,\n// in document order (the date paragraph, then each table-cell paragraph).\nconst replacements = [\n  \"2023-10-27 Friday\",\n  \"87\u00f78=\", \"63\u00f78=\", \"16\u00f77=\", \"97\u00f76=\", \"51\u00f73=\",\n  \"64\u00f72=\", \"32\u00f75=\", \"97\u00f79=\", \"45\u00f72=\", \"46\u00f75=\",\n  \"43\u00f74=\", \"90\u00f78=\", \"16\u00f76=\", \"19\u00f79=\", \"79\u00f72=\",\n  \"89\u00f74=\", \"49\u00f76=\", \"22\u00f78=\", \"11\u00f76=\", \"24\u00f77=\",\n  \"65\u00f77=\", \"53\u00f79=\", \"49\u00f78=\", \"13\u00f78=\", \"73\u00f72=\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet i = 0;\nfor (const paragraph of paragraphs.items) {\n  const current = paragraph.text;\n  if (current && current.length > 0) {\n    if (i < replacements.length) {\n      paragraph.insertText(replacements[i], \"Replace\");\n    }\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date line at the top of the worksheet.\n$d.Paragraphs(1).Range.Text = \"2023-10-27 Friday\"\n\n# Update every \"a\u00f7b=\" problem in the practice table (5 populated rows x 5 columns,\n# with 3 blank rows between each populated row).\n$t = $d.Tables(1)\n\n$rowNumbers = @(1, 5, 9, 13, 17)\n$rowValues = @(\n  @(\"87\u00f78=\", \"63\u00f78=\", \"16\u00f77=\", \"97\u00f76=\", \"51\u00f73=\"),\n  @(\"64\u00f72=\", \"32\u00f75=\", \"97\u00f79=\", \"45\u00f72=\", \"46\u00f75=\"),\n  @(\"43\u00f74=\", \"90\u00f78=\", \"16\u00f76=\", \"19\u00f79=\", \"79\u00f72=\"),\n  @(\"89\u00f74=\", \"49\u00f76=\", \"22\u00f78=\", \"11\u00f76=\", \"24\u00f77=\"),\n  @(\"65\u00f77=\", \"53\u00f79=\", \"49\u00f78=\", \"13\u00f78=\", \"73\u00f72=\")\n)\n\nfor ($i = 0; $i -lt $rowNumbers.Count; $i++) {\n  $row = $rowNumbers[$i]\n  $values = $rowValues[$i]\n  for ($col = 1; $col -le 5; $col++) {\n    $t.Cell($row, $col).Range.Text = $values[$col - 1]\n  }\n}\n"}
